# Generate Report for Handoff
#
# A fresh handoff run produced a new GUID-named source markdown file and a
# new handoff package (.xlf) for each locale, with fresh handoff
# timestamps. Update every cell that references the old identifiers /
# timestamps to the new ones (the "File Name" column on each locale sheet,
# the handoff file name, and the handoff datetime).

$wb = $excel.ActiveWorkbook

$oldGuid = "3a5cee9c-f299-4390-afe5-cd9f6bb7d837"
$newGuid = "2f8cc5e9-36bd-4329-8fa0-75d10cbea17f"

$oldHash = "d445560a492c5e341f1f8866ac8f67fdd62e1c1a"
$newHash = "21f0e2e366d655d60e21386f30c24a2104a48801"

$newMdName = "$newGuid.md"

$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newZhHandoffTime = "2016-03-09 10:02:25"

$newDeXlfName = "$newGuid.$newHash.de-de.xlf"
$newDeHandoffTime = "2016-03-09 10:02:34"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Source .md file name got a new GUID - reflected in the "File Name" / first
# column of every sheet.
$wsOverview.Range("A2").Value = $newMdName
$wsZh.Range("A2").Value = $newMdName
$wsDe.Range("A2").Value = $newMdName

# zh-cn: new handoff package name + new handoff datetime.
$wsZh.Range("C2").Value = $newZhXlfName
$wsZh.Range("D2").Value = $newZhHandoffTime

# de-de: new handoff package name + new handoff datetime.
$wsDe.Range("C2").Value = $newDeXlfName
$wsDe.Range("D2").Value = $newDeHandoffTime
